# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gets a new "property_category" column inserted
# between "total" and "date" (so date/legislator_name/legislator_id each
# shift one column to the right, H->I, I->J, J->K), populated with the
# literal value "stock" for every data row. Also fix a stray embedded
# space in one of the company names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H, pushing date/legislator_name/legislator_id
# (and the dimension/spans) one column to the right; formatting (styles)
# of the shifted cells is carried along automatically.
$ws.Range("H1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("H1").Value = "property_category"

# Populate the new column for each data row with the property category.
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"

# Fix the stray space in the company name on row 4 (index 92).
$ws.Range("B4").Value = "宸鴻光電科技股份有限公司"
